$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old Row 7 values
$ws.Range("D2").Value = 44742
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 100
$ws.Range("P2").Value = 14500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 806

# Row 3 <- old Row 5 values
$ws.Range("D3").Value = 44330
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("S3").Value = 861
$ws.Range("T3").Value = 18

# Row 4 <- old Row 2 values
$ws.Range("D4").Value = 44719
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14400
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región del Maule"
$ws.Range("S4").Value = 800
$ws.Range("T4").Value = 18

# Row 5 <- old Row 4 values
$ws.Range("D5").Value = 44707
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = "$/caja 12 kilos empedrada"
$ws.Range("S5").Value = 1042
$ws.Range("T5").Value = 12

# Row 7 <- old Row 3 values
$ws.Range("D7").Value = 44708
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12571
$ws.Range("Q7").Value = "$/caja 12 kilos empedrada"
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 1048
$ws.Range("T7").Value = 12
